# Applies the "Updated cryptos list" data refresh to the Price (D) and
# Volume(1h) (E) columns, plus the row 51 coin swap (ThetaToken -> Cronos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.081.91"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "'3.421.78"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'409.37"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'129.89"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("D7").Value = "'0.632"
$ws.Range("E7").Value = "  +6.72%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +10.85%  "
$ws.Range("D10").Value = "'0.150"
$ws.Range("E10").Value = "  +22.90%  "
$ws.Range("D11").Value = "'42.85"
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("D12").Value = "'0.0000230"
$ws.Range("E12").Value = "  +74.08%  "
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("D14").Value = "'3.968.80"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("E15").Value = "  +6.37%  "
$ws.Range("D16").Value = "'21.11"
$ws.Range("E16").Value = "  +5.84%  "
$ws.Range("D17").Value = "'3.413.23"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "'12.38"
$ws.Range("E18").Value = "  +12.61%  "
$ws.Range("E19").Value = "  +4.47%  "
$ws.Range("D20").Value = "'62.011.52"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'401.29"
$ws.Range("E21").Value = "  +26.76%  "
$ws.Range("D22").Value = "'89.86"
$ws.Range("E22").Value = "  +6.13%  "
$ws.Range("E23").Value = "  -0.96%  "
$ws.Range("D24").Value = "'13.18"
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("E25").Value = "  +3.75%  "
$ws.Range("D26").Value = "'32.82"
$ws.Range("E26").Value = "  +10.10%  "
$ws.Range("D27").Value = "'8.74"
$ws.Range("E27").Value = "  +5.92%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").Value = "'7.61"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("E31").Value = "  +0.18%  "
$ws.Range("E32").Value = "  -0.44%  "
$ws.Range("D33").Value = "'11.94"
$ws.Range("E33").Value = "  +4.67%  "
$ws.Range("D34").Value = "'43.04"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").Value = "  +0.66%  "
$ws.Range("D36").Value = "'0.0498"
$ws.Range("E36").Value = "  +2.88%  "
$ws.Range("D37").Value = "'54.14"
$ws.Range("E37").Value = "  +4.33%  "
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").Value = "'0.134"
$ws.Range("E40").Value = "  +7.51%  "
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("E42").Value = "  +7.04%  "
$ws.Range("D43").Value = "'141.94"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("E44").Value = "  -1.09%  "
$ws.Range("D45").Value = "'4.10"
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("D46").Value = "'2.42"
$ws.Range("E46").Value = "  +9.15%  "
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").Value = "'21.69"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("D49").Value = "'2.122.21"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").Value = "'2.37"
$ws.Range("E50").Value = "  +2.63%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.131"
$ws.Range("E51").Value = "  +16.59%  "
